# Second draft of quarterly report: update the 'Revised' figures in the
# cycle-by-cycle summary table (Revised, n_percent_zero_Revised,
# n_percent_seven_plus_Revised columns) to match James's corrected numbers.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$mismatchCount = 0

function Replace-CellText($table, $row, $col, $old, $new) {
    # Build a fresh Range from the cell's Start/End so Find.Execute stays
    # strictly scoped to this one cell and can't bleed into other cells
    # that happen to contain the same substring (e.g. "4.5" inside "74.5").
    $cell = $table.Cell($row, $col)
    $rng = $d.Range($cell.Range.Start, $cell.Range.End)
    $before = $rng.Text
    $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 0, $false, $new, 1) | Out-Null

    $cellAfter = $table.Cell($row, $col)
    $after = $cellAfter.Range.Text
    if ($after -notmatch [regex]::Escape($new)) {
        Write-Host "WARNING: Cell($row,$col) expected to contain '$new' but got '$after' (was '$before')"
        $script:mismatchCount++
    }
}

Replace-CellText $t 2 5 "11 (11.2)" "11 ( 9.8)"
Replace-CellText $t 2 7 "4 (4.1)" "6 (5.4)"
Replace-CellText $t 3 3 "8.0" "8.5"
Replace-CellText $t 3 5 "22 (22.4)" "24 (21.4)"
Replace-CellText $t 3 7 "0 (0.0)" "1 (0.9)"
Replace-CellText $t 4 3 "5.5" "5.0"
Replace-CellText $t 4 5 "34 (34.7)" "37 (33.0)"
Replace-CellText $t 4 7 "2 (2.0)" "2 (1.8)"
Replace-CellText $t 5 3 "4.5" "4.0"
Replace-CellText $t 5 5 "41 (41.8)" "48 (42.9)"
Replace-CellText $t 5 7 "4 (4.1)" "5 (4.5)"
Replace-CellText $t 6 5 "42 (42.9)" "49 (43.8)"
Replace-CellText $t 6 7 "1 (1.0)" "1 (0.9)"
Replace-CellText $t 7 3 "2.5" "2.0"
Replace-CellText $t 7 5 "47 (48.0)" "55 (49.1)"
Replace-CellText $t 7 7 "3 (3.1)" "3 (2.7)"
Replace-CellText $t 8 5 "56 (57.1)" "63 (56.2)"
Replace-CellText $t 8 7 "0 (0.0)" "2 (1.8)"
Replace-CellText $t 9 5 "53 (54.1)" "60 (53.6)"
Replace-CellText $t 9 7 "1 (1.0)" "1 (0.9)"
Replace-CellText $t 10 5 "58 (59.2)" "65 (58.0)"
Replace-CellText $t 10 7 "2 (2.0)" "2 (1.8)"
Replace-CellText $t 11 5 "60 (61.2)" "68 (60.7)"
Replace-CellText $t 11 7 "1 (1.0)" "1 (0.9)"
Replace-CellText $t 12 5 "62 (63.3)" "70 (62.5)"
Replace-CellText $t 13 5 "65 (66.3)" "74 (66.1)"
Replace-CellText $t 13 7 "2 (2.0)" "2 (1.8)"

if ($mismatchCount -eq 0) {
    Write-Host "All 27 revised-figure cells updated successfully."
} else {
    Write-Host "$mismatchCount cell(s) did not update as expected."
}
